# Rename several pf_* result-column headers across all fault-type sheets.
# "deg" -> "degree", and the mis-named reactive-power / voltage-magnitude
# columns get corrected (mw -> mvar, *_bus_pu -> *_pu, *_bus_deg -> *_degree).

$wb = $excel.ActiveWorkbook

# --- Sheets 1-8 (LLL_*): 3-phase aggregate columns only ----------------
$sheets3ph = @(
    "LLL_max_6", "LLL_max_10", "LLL_max_fault_6", "LLL_max_fault_10",
    "LLL_min_6", "LLL_min_10", "LLL_min_fault_6", "LLL_min_fault_10"
)

$headerMap3ph = @{
    "L1" = "pf_ikss_from_degree"
    "M1" = "pf_ikss_to_degree"
    "P1" = "pf_va_from_degree"
    "Q1" = "pf_va_to_degree"
}

foreach ($name in $sheets3ph) {
    $ws = $wb.Worksheets.Item($name)
    foreach ($addr in $headerMap3ph.Keys) {
        $ws.Range($addr).Value = $headerMap3ph[$addr]
    }
}

# --- Sheets 9-32 (LL_*, LLG_*, LG_*): per-phase a/b/c columns -----------
$sheetsPhase = @(
    "LL_max_6", "LL_max_10", "LL_max_fault_6", "LL_max_fault_10",
    "LL_min_6", "LL_min_10", "LL_min_fault_6", "LL_min_fault_10",
    "LLG_max_6", "LLG_max_10", "LLG_max_fault_6", "LLG_max_fault_10",
    "LLG_min_6", "LLG_min_10", "LLG_min_fault_6", "LLG_min_fault_10",
    "LG_max_6", "LG_max_10", "LG_max_fault_6", "LG_max_fault_10",
    "LG_min_6", "LG_min_10", "LG_min_fault_6", "LG_min_fault_10"
)

$headerMapPhase = @{
    "T1"  = "pf_q_a_from_mvar"
    "U1"  = "pf_q_b_from_mvar"
    "V1"  = "pf_q_c_from_mvar"
    "W1"  = "pf_q_a_to_mvar"
    "X1"  = "pf_q_b_to_mvar"
    "Y1"  = "pf_q_c_to_mvar"
    "Z1"  = "pf_ikss_a_from_degree"
    "AA1" = "pf_ikss_b_from_degree"
    "AB1" = "pf_ikss_c_from_degree"
    "AC1" = "pf_ikss_a_to_degree"
    "AD1" = "pf_ikss_b_to_degree"
    "AE1" = "pf_ikss_c_to_degree"
    "AG1" = "pf_vm_b_from_pu"
    "AH1" = "pf_vm_c_from_pu"
    "AI1" = "pf_vm_a_to_pu"
    "AJ1" = "pf_vm_b_to_pu"
    "AK1" = "pf_vm_c_to_pu"
    "AL1" = "pf_va_a_from_degree"
    "AM1" = "pf_va_b_from_degree"
    "AN1" = "pf_va_c_from_degree"
    "AO1" = "pf_va_a_to_degree"
    "AP1" = "pf_va_b_to_degree"
    "AQ1" = "pf_va_c_to_degree"
}

foreach ($name in $sheetsPhase) {
    $ws = $wb.Worksheets.Item($name)
    foreach ($addr in $headerMapPhase.Keys) {
        $ws.Range($addr).Value = $headerMapPhase[$addr]
    }
}
